$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# ---------------------------------------------------------------------------
# Append 4 new rows (102-105) to the TaskList sheet, continuing the existing
# "S.No / Date / Task / ... / Time Taken" log with the 11th-day entries.
# ---------------------------------------------------------------------------

# Row 102 (S.No 101) - style twin of row 98/100 (A = no-wrap border style)
$ws.Range("A98").Copy()
$ws.Range("A102").PasteSpecial(-4122)
$ws.Range("B101").Copy()
$ws.Range("B102").PasteSpecial(-4122)
$ws.Range("D101:E101").Copy()
$ws.Range("D102:E102").PasteSpecial(-4122)
$ws.Range("C101").Copy()
$ws.Range("C102").PasteSpecial(-4122)
$ws.Range("F101").Copy()
$ws.Range("F102").PasteSpecial(-4122)
$ws.Range("G101:H101").Copy()
$ws.Range("G102:H102").PasteSpecial(-4122)

$ws.Range("A102").Value = 101
$ws.Range("B102").Value = 42803
$ws.Range("C102").Value = "Webflow"
$ws.Range("F102").Value = "6hrs"

# Row 103 (S.No 102) - style twin of row 99/101 (A = wrap border style)
$ws.Range("A99").Copy()
$ws.Range("A103").PasteSpecial(-4122)
$ws.Range("B101").Copy()
$ws.Range("B103").PasteSpecial(-4122)
$ws.Range("D101:E101").Copy()
$ws.Range("D103:E103").PasteSpecial(-4122)
$ws.Range("C101").Copy()
$ws.Range("C103").PasteSpecial(-4122)
$ws.Range("F101").Copy()
$ws.Range("F103").PasteSpecial(-4122)
$ws.Range("G101:H101").Copy()
$ws.Range("G103:H103").PasteSpecial(-4122)

$ws.Range("A103").Value = 102
$ws.Range("B103").Value = 42803
$ws.Range("C103").Value = "Spring Security"
$ws.Range("F103").Value = "6hrs"

# Row 104 (S.No 103) - style twin of row 98/100 (A = no-wrap border style)
$ws.Range("A100").Copy()
$ws.Range("A104").PasteSpecial(-4122)
$ws.Range("B101").Copy()
$ws.Range("B104").PasteSpecial(-4122)
$ws.Range("D101:E101").Copy()
$ws.Range("D104:E104").PasteSpecial(-4122)
$ws.Range("C101").Copy()
$ws.Range("C104").PasteSpecial(-4122)
$ws.Range("F101").Copy()
$ws.Range("F104").PasteSpecial(-4122)
$ws.Range("G101:H101").Copy()
$ws.Range("G104:H104").PasteSpecial(-4122)

$ws.Range("A104").Value = 103
$ws.Range("B104").Value = 42803
$ws.Range("C104").Value = "Cart Implementation "
$ws.Range("F104").Value = "Pending"

# Row 105 (S.No 104) - style twin of row 99/101 (A = wrap border style)
$ws.Range("A101").Copy()
$ws.Range("A105").PasteSpecial(-4122)
$ws.Range("B101").Copy()
$ws.Range("B105").PasteSpecial(-4122)
$ws.Range("D101:E101").Copy()
$ws.Range("D105:E105").PasteSpecial(-4122)
$ws.Range("C101").Copy()
$ws.Range("C105").PasteSpecial(-4122)
$ws.Range("F101").Copy()
$ws.Range("F105").PasteSpecial(-4122)
$ws.Range("G101:H101").Copy()
$ws.Range("G105:H105").PasteSpecial(-4122)

$ws.Range("A105").Value = 104
$ws.Range("B105").Value = 42803
$ws.Range("C105").Value = "Documentation"
$ws.Range("F105").Value = "30Minutes"

# ---------------------------------------------------------------------------
# Update the view so the newly-added rows are visible / selected, matching
# the scroll position and active cell captured after today's entries.
# ---------------------------------------------------------------------------
$ws.Application.GoTo($ws.Range("C99"), $True)
$ws.Range("A85").Select()
$ws.Range("C99").Select()
